# ==========================================================
# Weekly CompStat report refresh: new crime data collected
# ==========================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume/issue number and week-covering dates ---
$ws.Range("A8").Value = "Volume 33   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/26/2026  Through  2/1/2026"

# --- Crime statistics table (rows 16-31): refreshed weekly figures ---
# --- Row 16 ---
$c = $ws.Range("C16"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("E16"); $c.Value = 0
$c = $ws.Range("F16"); $c.Value = 3
$c = $ws.Range("G16"); $c.Value = 4
$c = $ws.Range("H16"); $c.Value = -25
$c = $ws.Range("I16"); $c.Value = 4
$c = $ws.Range("J16"); $c.Value = 4
$c = $ws.Range("L16"); $c.Value = -60
$c = $ws.Range("M16"); $c.Value = -77.777777777777
$c = $ws.Range("N16"); $c.Value = -95.238095238095

# --- Row 17 ---
$c = $ws.Range("C17"); $c.Value = 5
$c = $ws.Range("D17"); $c.Value = 3
$c = $ws.Range("E17"); $c.Value = 66.666666666666
$c = $ws.Range("F17"); $c.Value = 13
$c = $ws.Range("G17"); $c.Value = 14
$c = $ws.Range("H17"); $c.Value = -7.142857142857
$c = $ws.Range("I17"); $c.Value = 13
$c = $ws.Range("J17"); $c.Value = 17
$c = $ws.Range("K17"); $c.Value = -23.529411764705
$c = $ws.Range("L17"); $c.Value = 44.444444444444
$c = $ws.Range("M17"); $c.Value = 333.333333333333
$c = $ws.Range("N17"); $c.Value = -27.777777777777

# --- Row 18 ---
$c = $ws.Range("D18"); $c.Value = 4
$c = $ws.Range("F18"); $c.Value = 5
$c = $ws.Range("G18"); $c.Value = 14
$c = $ws.Range("H18"); $c.Value = -64.285714285714
$c = $ws.Range("J18"); $c.Value = 15
$c = $ws.Range("K18"); $c.Value = -53.333333333333
$c = $ws.Range("L18"); $c.Value = -50
$c = $ws.Range("M18"); $c.Value = -46.153846153846
$c = $ws.Range("N18"); $c.Value = -91.25

# --- Row 19 ---
$c = $ws.Range("C19"); $c.Value = 6
$c = $ws.Range("E19"); $c.Value = 50
$c = $ws.Range("F19"); $c.Value = 17
$c = $ws.Range("G19"); $c.Value = 22
$c = $ws.Range("H19"); $c.Value = -22.727272727272
$c = $ws.Range("I19"); $c.Value = 20
$c = $ws.Range("J19"); $c.Value = 26
$c = $ws.Range("K19"); $c.Value = -23.076923076923
$c = $ws.Range("L19"); $c.Value = -42.857142857142
$c = $ws.Range("M19"); $c.Value = -47.368421052631
$c = $ws.Range("N19"); $c.Value = -51.219512195122

# --- Row 20 ---
$c = $ws.Range("C20"); $c.NumberFormat = "@"; $c.Value = "0"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("D20"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("E20"); $c.Value = -100
$ws.Range("L14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("J20"); $c.Value = 2
$c = $ws.Range("K20"); $c.Value = 50
$c = $ws.Range("L20"); $c.Value = -70
$c = $ws.Range("M20"); $c.Value = -50
$c = $ws.Range("N20"); $c.Value = -97.247706422018

# --- Row 21 ---
$c = $ws.Range("C21"); $c.Value = 12
$c = $ws.Range("D21"); $c.Value = 13
$c = $ws.Range("E21"); $c.Value = -7.692307692307
$c = $ws.Range("G21"); $c.Value = 56
$c = $ws.Range("H21"); $c.Value = -26.785714285714
$c = $ws.Range("I21"); $c.Value = 47
$c = $ws.Range("J21"); $c.Value = 65
$c = $ws.Range("K21"); $c.Value = -27.692307692307
$c = $ws.Range("L21"); $c.Value = -41.25
$c = $ws.Range("M21"); $c.Value = -39.743589743589
$c = $ws.Range("N21"); $c.Value = -86.011904761904

# --- Row 22 ---
$c = $ws.Range("D22"); $c.NumberFormat = "@"; $c.Value = "0"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("E22"); $c.NumberFormat = "@"; $c.Value = "***.*"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("F22"); $c.NumberFormat = "@"; $c.Value = "0"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("H22"); $c.Value = -100
$c = $ws.Range("L22"); $c.Value = -66.666666666666
$c = $ws.Range("M22"); $c.Value = -80

# --- Row 23 ---
$c = $ws.Range("F23"); $c.NumberFormat = "@"; $c.Value = "0"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("G23"); $c.NumberFormat = "@"; $c.Value = "0"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("H23"); $c.NumberFormat = "@"; $c.Value = "***.*"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false

# --- Row 24 ---
$c = $ws.Range("C24"); $c.Value = 24
$c = $ws.Range("D24"); $c.Value = 17
$c = $ws.Range("E24"); $c.Value = 41.176470588235
$c = $ws.Range("F24"); $c.Value = 104
$c = $ws.Range("G24"); $c.Value = 100
$c = $ws.Range("H24"); $c.Value = 4
$c = $ws.Range("I24"); $c.Value = 112
$c = $ws.Range("J24"); $c.Value = 108
$c = $ws.Range("K24"); $c.Value = 3.703703703703
$c = $ws.Range("L24"); $c.Value = 19.148936170212
$c = $ws.Range("M24"); $c.Value = 55.555555555555

# --- Row 25 ---
$c = $ws.Range("C25"); $c.Value = 17
$c = $ws.Range("D25"); $c.Value = 13
$c = $ws.Range("E25"); $c.Value = 30.769230769230
$c = $ws.Range("F25"); $c.Value = 68
$c = $ws.Range("G25"); $c.Value = 67
$c = $ws.Range("H25"); $c.Value = 1.492537313432
$c = $ws.Range("I25"); $c.Value = 75
$c = $ws.Range("J25"); $c.Value = 71
$c = $ws.Range("K25"); $c.Value = 5.633802816901
$c = $ws.Range("L25"); $c.Value = 31.578947368421

# --- Row 26 ---
$c = $ws.Range("C26"); $c.Value = 5
$c = $ws.Range("D26"); $c.Value = 1
$c = $ws.Range("E26"); $c.Value = 400
$c = $ws.Range("F26"); $c.Value = 12
$c = $ws.Range("G26"); $c.Value = 19
$c = $ws.Range("H26"); $c.Value = -36.842105263157
$c = $ws.Range("I26"); $c.Value = 12
$c = $ws.Range("J26"); $c.Value = 24
$c = $ws.Range("K26"); $c.Value = -50
$c = $ws.Range("L26"); $c.Value = -7.692307692307
$c = $ws.Range("M26"); $c.Value = -40

# --- Row 27 ---
$c = $ws.Range("C27"); $c.NumberFormat = "@"; $c.Value = "0"
$ws.Range("C14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false

# --- Row 28 ---
$c = $ws.Range("D28"); $c.Value = 2
$c = $ws.Range("G28"); $c.Value = 5
$c = $ws.Range("H28"); $c.Value = -40
$c = $ws.Range("J28"); $c.Value = 5
$c = $ws.Range("K28"); $c.Value = -40
$c = $ws.Range("L28"); $c.Value = -75

# --- Row 29 ---
$c = $ws.Range("C29"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("F29"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("I29"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("L29"); $c.Value = 0
$c = $ws.Range("N29"); $c.Value = 0
$ws.Range("L14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false

# --- Row 30 ---
$c = $ws.Range("C30"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("F30"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("I30"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("L30"); $c.Value = 0
$c = $ws.Range("N30"); $c.Value = 0
$ws.Range("L14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false

# --- Row 31 ---
$c = $ws.Range("D31"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("E31"); $c.Value = -100
$ws.Range("L14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("G31"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("H31"); $c.Value = -100
$ws.Range("L14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("J31"); $c.Value = 1
$ws.Range("G15").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
$c = $ws.Range("K31"); $c.Value = -100
$ws.Range("L14").Copy(); $c.PasteSpecial(-4122); $excel.CutCopyMode = $false
